# Apply the "Cập nhật thêm thông tin vài mẫu thiết kế." commit.
#
# Summary of changes:
#  - Creational sheet: fill in notes for Builder / Prototype / Singleton rows.
#  - Behavioral sheet: fill in notes for Chain of Responsibility / Command /
#    Interpreter / Iterator / Mediator rows, widen column C, move the
#    selection, and make this the active tab of the workbook.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Creational sheet
# ---------------------------------------------------------------------------
$creational = $wb.Worksheets.Item("Creational")

$creational.Range("C7").Value = "Việc khởi tạo thực hiện duy nhất 1 lần"
$creational.Range("D5").Value = "Chưa rõ"
$creational.Range("D6").Value = "Chưa rõ"
$creational.Range("D7").Value = "Có thể áp dụng"

# ---------------------------------------------------------------------------
# Behavioral sheet
# ---------------------------------------------------------------------------
$behavioral = $wb.Worksheets.Item("Behavioral")

$behavioral.Range("C3").Value = "Quyết định quy trình chạy của các lớp theo thứ tự nhất định"
$behavioral.Range("D3").Value = "Chưa rõ"

$behavioral.Range("C4").Value = "Chuyên xử lý các hành động undo, redo"
$behavioral.Range("D4").Value = "Chưa rõ"

$behavioral.Range("C5").Value = "Quản lý định dạng date hoặc việc đọc dữ liệu từ bên ngoài của nhiều class"
$behavioral.Range("D5").Value = "Chưa rõ"

$behavioral.Range("C6").Value = "Quản lý kiểu danh sách nhiều phần tử"
$behavioral.Range("D6").Value = "Chưa rõ"

$behavioral.Range("C7").Value = "Điều phối thông điệp với các thể hiện khác nhau"
$behavioral.Range("D7").Value = "Chưa rõ"

$behavioral.Columns.Item(3).ColumnWidth = 62

# Move the selection / make Behavioral the active sheet + tab.
$behavioral.Range("C7").Select()
$behavioral.Activate()
